$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# 1) Update the lookup table on Sheet3 (B20:B36) with new values.
#    Formulas on Sheet3 (C2:C18) and on Sheet1 (CB/CC columns) reference
#    this table via VLOOKUP and will recalculate automatically.
$ws3.Range("B20").Value = 12.739989309477449   # 3D QUESO 92GX27
$ws3.Range("B21").Value = 3.1869787649075585   # CHEETOS QUESO 85GX24X1
$ws3.Range("B22").Value = 8.6412920283286159   # DORITOS QUESO 129GX19
$ws3.Range("B23").Value = 7.9244268301531733   # DORITOS QUESO 70X40G
$ws3.Range("B24").Value = 9.0492200316057527   # DORITOS QUESO 77GX26
$ws3.Range("B25").Value = 0                    # LAYS CEBOLLA CARAMELIZADA 85GX25
$ws3.Range("B26").Value = 15.857791621052723   # LAYS CLASICAS 145GRX18
$ws3.Range("B27").Value = 2.4598760302145206   # LAYS CLASICAS 249GRX14
$ws3.Range("B28").Value = 6.5045164785246206   # LAYS CLASICAS 40GX68
$ws3.Range("B29").Value = 6.0290405848472108   # LAYS CLASICAS 85GX25
$ws3.Range("B30").Value = 3.3922800000104623   # LAYS ONDAS FH 30GX72
$ws3.Range("B31").Value = 5.0789970000861651   # LAYS ONDAS FH 70GX28
$ws3.Range("B32").Value = 19.219728041853248   # LAYS QSO Y CEBOLLA 34GX72
$ws3.Range("B33").Value = 3.8432090075392069   # PEHUAMAR ACANALADA 520GX9
$ws3.Range("B34").Value = 10.064784158902867   # PEHUAMAR MAICITOS 285GX10
$ws3.Range("B35").Value = 8.3579217224649298   # PEHUAMAR PAPA LISA 520GX9
$ws3.Range("B36").Value = 59.091037666315096   # QUAKER AVENA INSTANT FORTIF 18X280G

$excel.Calculate()

# 2) Add a new "27-nov" snapshot column (CQ) on Sheet1, right after CP,
#    holding the freshly recalculated VLOOKUP values as static numbers.
$ws1.Range("CQ1").Value = "27-nov"

for ($r = 2; $r -le 18; $r++) {
    $src = $ws1.Cells.Item($r, 80)   # column CB = 80th column, holds the recalculated VLOOKUP value
    $dst = $ws1.Cells.Item($r, 95)   # column CQ = 95th column
    $dst.Value = $src.Value2
    $dst.NumberFormat = $src.NumberFormat
}

# 3) Update the sheet view: drop the frozen/offset top-left cell and move
#    the active selection.
$ws1.Range("CK29").Select()
